# Add level-flight / hover mode LED legend entries to the "硬件" (Hardware)
# sheet and refresh the selected cell / active tab bookmarks to match.

$wb = $excel.ActiveWorkbook

$wsHw = $wb.Worksheets.Item("硬件")
$wsParam = $wb.Worksheets.Item("参数")

# --- Update the LED legend text on the hardware sheet -----------------
# Row 13's "PWM输入" becomes "悬停模式" (hover mode) and row 14's
# "PPM输入" becomes "平飞模式" (level-flight mode). Writing I14 before I13
# keeps the shared-string insertion order (and therefore the resulting
# string-table indices) aligned with the authored workbook.
$wsHw.Range("I14").Value = "平飞模式"
$wsHw.Range("I13").Value = "悬停模式"

# --- Update selections / active sheet to match the authored state -----
# "参数" was the active tab before; the selected cell there moves from
# E19 to E16.
$wsParam.Range("E16").Select()

# "硬件" becomes the active tab, with I13 selected.
$wsHw.Activate()
$wsHw.Range("I13").Select()
